$d = $word.ActiveDocument

$replacements = @(
    @{ old = "2024-09-25 Wednesday"; new = "2024-09-26 Thursday" },
    @{ old = "252÷4=63, 0";  new = "632÷4=158, 0" },
    @{ old = "710÷4=177, 2"; new = "184÷7=26, 2" },
    @{ old = "662÷3=220, 2"; new = "817÷7=116, 5" },
    @{ old = "355÷2=177, 1"; new = "958÷6=159, 4" },
    @{ old = "303÷4=75, 3";  new = "912÷6=152, 0" },
    @{ old = "581÷4=145, 1"; new = "481÷8=60, 1" },
    @{ old = "395÷5=79, 0";  new = "838÷2=419, 0" },
    @{ old = "770÷2=385, 0"; new = "262÷5=52, 2" },
    @{ old = "335÷4=83, 3";  new = "715÷6=119, 1" },
    @{ old = "842÷9=93, 5";  new = "595÷5=119, 0" },
    @{ old = "592÷7=84, 4";  new = "457÷6=76, 1" },
    @{ old = "925÷7=132, 1"; new = "966÷9=107, 3" },
    @{ old = "157÷5=31, 2";  new = "728÷3=242, 2" },
    @{ old = "841÷2=420, 1"; new = "555÷3=185, 0" },
    @{ old = "464÷9=51, 5";  new = "637÷3=212, 1" },
    @{ old = "185÷2=92, 1";  new = "993÷3=331, 0" },
    @{ old = "180÷9=20, 0";  new = "974÷9=108, 2" },
    @{ old = "977÷5=195, 2"; new = "389÷4=97, 1" },
    @{ old = "239÷8=29, 7";  new = "494÷6=82, 2" },
    @{ old = "653÷6=108, 5"; new = "668÷8=83, 4" },
    @{ old = "294÷4=73, 2";  new = "451÷9=50, 1" },
    @{ old = "136÷2=68, 0";  new = "871÷9=96, 7" },
    @{ old = "251÷9=27, 8";  new = "708÷8=88, 4" },
    @{ old = "896÷8=112, 0"; new = "268÷5=53, 3" },
    @{ old = "152÷9=16, 8";  new = "212÷4=53, 0" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
